$wb = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item(1)   # "orders" sheet
$ws2 = $wb.Worksheets.Item(2)   # "customer_info" sheet

# orders sheet: the old customer-detail columns (Email, Phone, Name,
# Subteam) are replaced with the new merch item columns, and a new
# "Price" column is appended, so the header row becomes:
# Order Id | Item Id | Size | Quantity | Colour | Price
$ws1.Range("B1").Value = "Item Id"
$ws1.Range("C1").Value = "Size"
$ws1.Range("D1").Value = "Quantity"
$ws1.Range("E1").Value = "Colour"
$ws1.Range("F1").Value = "Price"

# customer_info sheet: keep the existing columns and append a new
# "Additional Notes" field at the end (column Q).
$ws2.Range("Q1").Value = "Additional Notes"
